$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price / volume data as per the scheduled GitHub Actions refresh.
# For cells whose new value looks like a plain number (e.g. "1.00", "6.97"), force a Text
# number format first so Excel keeps the original textual representation instead of
# silently converting it to a numeric value (which would drop formatting / precision).

$ws.Range("D2").Value = '67.418.43'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").Value = '3.441.26'
$ws.Range("E3").Value = '  -2.45%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.68'
$ws.Range("E5").Value = '  -1.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.37'
$ws.Range("E6").Value = '  -3.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.604'
$ws.Range("E8").Value = '  +1.44%  '
$ws.Range("D9").Value = '3.440.74'
$ws.Range("E9").Value = '  -2.45%  '
$ws.Range("E10").Value = '  -1.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.97'
$ws.Range("E11").Value = '  -2.63%  '
$ws.Range("E12").Value = '  -4.49%  '
$ws.Range("D13").Value = '4.044.92'
$ws.Range("E13").Value = '  -2.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '31.85'
$ws.Range("E14").Value = '  -2.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.132'
$ws.Range("E15").Value = '  -1.39%  '
$ws.Range("D16").Value = '67.452.42'
$ws.Range("E16").Value = '  -0.80%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000175'
$ws.Range("E17").Value = '  -3.61%  '
$ws.Range("D18").Value = '3.444.06'
$ws.Range("E18").Value = '  -2.50%  '
$ws.Range("E19").Value = '  -4.94%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.95'
$ws.Range("E20").Value = '  -6.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '385.07'
$ws.Range("E21").Value = '  -3.85%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.81'
$ws.Range("E22").Value = '  -3.92%  '
$ws.Range("E23").Value = '  +2.02%  '
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("E25").Value = '  -3.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '71.13'
$ws.Range("E26").Value = '  -3.80%  '
$ws.Range("E27").Value = '  -5.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.17'
$ws.Range("E28").Value = '  -4.45%  '
$ws.Range("E29").Value = '  -2.35%  '
$ws.Range("E30").Value = '  +0.29%  '
$ws.Range("E31").Value = '  -4.37%  '
$ws.Range("E33").Value = '  -7.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.42'
$ws.Range("E34").Value = '  -3.19%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  -3.76%  '
$ws.Range("E37").Value = '  -7.82%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '161.17'
$ws.Range("E38").Value = '  -1.73%  '
$ws.Range("E39").Value = '  -0.41%  '
$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.72'
$ws.Range("E40").Value = '  -3.39%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.84'
$ws.Range("E41").Value = '  -5.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.60'
$ws.Range("E42").Value = '  -7.79%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.51'
$ws.Range("E43").Value = '  -5.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '25.74'
$ws.Range("E44").Value = '  -5.89%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0710'
$ws.Range("E45").Value = '  -4.58%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.85'
$ws.Range("E46").Value = '  -6.24%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.692.90'
$ws.Range("E47").Value = '  -7.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '41.14'
$ws.Range("E48").Value = '  -3.40%  '
$ws.Range("E49").Value = '  -3.78%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '325.24'
$ws.Range("E50").Value = '  -7.65%  '
$ws.Range("E51").Value = '  -5.28%  '
